$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused blank rows (old rows 41-52 become the new gap,
# collapsing the old row 53.. block up to row 43..) so that the former
# rows 53,55-64 become rows 43,45-54.
$ws.Rows("41:50").Delete() | Out-Null

# Rename the generic mapping name setting value.
$ws.Range("C46").Value = "Import depuis tableur (XLSX)"

# Bold the header rows (column titles + the two "Rule"/"Setting" legends).
$ws.Range("A1:K1").Font.Bold = $true
$ws.Range("B43:E43").Font.Bold = $true
$ws.Range("B45:E45").Font.Bold = $true

# Minor column width tweaks.
$ws.Columns("B").ColumnWidth = 8.75
$ws.Columns("C").ColumnWidth = 43.66
$ws.Columns("E").ColumnWidth = 31.75
$ws.Columns("F").ColumnWidth = 17.66
$ws.Columns("G").ColumnWidth = 27.75

# Page setup.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the last-saved selection.
$ws.Range("C8").Select() | Out-Null
